$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("heterogenous")

# New "MPS" section header (row 12), bold like the row 1 header
$ws.Range("A12").Value = "MPS"
$ws.Range("A12").Font.Bold = $true

# Fill the new parameter rows. Order chosen to reproduce shared-string table order.
$ws.Range("A16").Value = "TI_rivlen"

$ws.Range("A14").Value = "TI_rivwidth"
$ws.Range("B14").Value = "500m-2km"

$ws.Range("A15").Value = "TI_sinuosity"
$ws.Range("D15").Value = "27km of Salinas Valley = 37km of river"

$ws.Range("D13").Value = "Quick measurements from N Sal Valley"
$ws.Range("B13").Value = "2-4km"
$ws.Range("A13").Value = "TI_rivwavelength"

$ws.Range("B16").Value = "continuous"

$ws.Range("D14").Value = "Quick measurements from N Sal Valley"
$ws.Range("B15").Value = 1.341

# Update selection to mirror author's resulting cursor position
$ws.Range("A18").Select()
